$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
}

Set-CellText "D2" "65.066.23"
Set-CellText "E2" "  +0.74%  "
Set-CellText "D3" "3.373.73"
Set-CellText "E3" "  +0.35%  "
Set-CellText "E4" "  +0.00%  "
Set-CellText "D5" "555.10"
Set-CellText "E5" "  -0.07%  "
Set-CellText "D6" "174.03"
Set-CellText "E6" "  -1.02%  "
Set-CellText "D7" "0.630"
Set-CellText "E7" "  +1.90%  "
Set-CellText "D8" "3.359.72"
Set-CellText "E8" "  +0.21%  "
Set-CellText "E10" "  +5.72%  "
Set-CellText "D11" "0.635"
Set-CellText "E11" "  +0.83%  "
Set-CellText "D12" "53.57"
Set-CellText "E12" "  -1.62%  "
Set-CellText "E13" "  +1.39%  "
Set-CellText "D14" "9.15"
Set-CellText "E14" "  +0.53%  "
Set-CellText "D15" "3.917.06"
Set-CellText "E15" "  +0.54%  "
Set-CellText "D16" "18.32"
Set-CellText "E16" "  -0.63%  "
Set-CellText "D17" "0.119"
Set-CellText "E17" "  +0.45%  "
Set-CellText "D18" "3.371.01"
Set-CellText "E18" "  +0.37%  "
Set-CellText "D19" "64.997.85"
Set-CellText "E19" "  +0.75%  "
Set-CellText "D20" "11.83"
Set-CellText "E20" "  -0.16%  "
Set-CellText "E21" "  +0.92%  "
Set-CellText "D22" "458.86"
Set-CellText "E22" "  -1.00%  "
Set-CellText "D23" "4.88"
Set-CellText "E23" "  +1.86%  "
Set-CellText "E24" "  -0.35%  "
Set-CellText "E25" "  +5.96%  "
Set-CellText "D26" "87.82"
Set-CellText "E26" "  +1.92%  "
Set-CellText "D27" "2.89"
Set-CellText "E27" "  +1.37%  "
Set-CellText "D28" "10.66"
Set-CellText "E28" "  -2.53%  "
Set-CellText "E29" "  -1.18%  "
Set-CellText "D30" "31.02"
Set-CellText "E30" "  +2.94%  "
Set-CellText "E31" "  -1.87%  "
Set-CellText "D32" "63.23"
Set-CellText "E32" "  +7.41%  "
Set-CellText "D33" "11.43"
Set-CellText "E33" "  -0.44%  "
Set-CellText "D34" "576.87"
Set-CellText "E34" "  -0.89%  "
Set-CellText "E35" "  -0.68%  "
Set-CellText "E36" "  -0.08%  "
Set-CellText "D37" "3.65"
Set-CellText "E37" "  +4.23%  "
Set-CellText "E38" "  +1.48%  "
Set-CellText "E40" "  -1.38%  "
Set-CellText "D41" "0.0₃0738"
Set-CellText "E41" "  -2.78%  "
Set-CellText "D42" "3.094.49"
Set-CellText "E42" "  -0.28%  "
Set-CellText "E43" "  +1.04%  "
Set-CellText "E44" "  -1.90%  "
Set-CellText "E45" "  -0.40%  "
Set-CellText "B46" "Fetch.AI"
Set-CellText "C46" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-CellText "D46" "2.44"
Set-CellText "E46" "  -3.65%  "
Set-CellText "B47" "Stellar"
Set-CellText "C47" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-CellText "D47" "0.133"
Set-CellText "E47" "  +1.94%  "
Set-CellText "E48" "  +0.17%  "
Set-CellText "D49" "140.46"
Set-CellText "E49" "  +3.76%  "
Set-CellText "E50" "  -2.23%  "
Set-CellText "D51" "8.29"
Set-CellText "E51" "  -1.26%  "
